$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly answered row (row 11) of the "Tablo1" table with the
# ignore-file related Q&A that was added.
$ws.Range("F11").Value = "Yok"
$ws.Range("G11").Value = "Yok"
$ws.Range("H11").Value = "What is the difference between array & pointer array"
$ws.Range("I11").Value = "yok"

# Move the active selection to I11, matching where the author left the
# cursor after typing the new data.
$ws.Activate() | Out-Null
$ws.Range("I11").Select() | Out-Null
